$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.377.84"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "1.655.24"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "1.889.71"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "1.658.81"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("E15").Value = "  +3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "27.379.03"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.01%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "1.419.28"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.908"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("D46").Value = "1.797.90"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.04%  "
